$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell (outside the A1:E51 table) used to force text-typed numeric-looking
# values (e.g. "0.523") into the target cells without altering their cell style,
# since a direct Range.Value assignment of a plain numeric string would be
# auto-converted to a real number by Excel.
$scratch = $ws.Range('Z1')
$scratch.NumberFormat = '@'

$ws.Range('D2').Value = '29.492.39'
$ws.Range('E2').Value = '  +2.93%  '
$ws.Range('D3').Value = '1.608.74'
$ws.Range('E3').Value = '  +2.99%  '
$ws.Range('E4').Value = '  +0.20%  '
$scratch.Value = '212.91'
$scratch.Copy()
$ws.Range('D5').PasteSpecial(-4163)
$ws.Range('E5').Value = '  +1.26%  '
$scratch.Value = '0.523'
$scratch.Copy()
$ws.Range('D6').PasteSpecial(-4163)
$ws.Range('E6').Value = '  +7.13%  '
$ws.Range('E7').Value = '  +0.26%  '
$scratch.Value = '26.80'
$scratch.Copy()
$ws.Range('D8').PasteSpecial(-4163)
$ws.Range('E8').Value = '  +6.67%  '
$scratch.Value = '43.63'
$scratch.Copy()
$ws.Range('D9').PasteSpecial(-4163)
$ws.Range('E9').Value = '  -0.73%  '
$ws.Range('E10').Value = '  +2.82%  '
$ws.Range('E11').Value = '  +2.89%  '
$ws.Range('E12').Value = '  +1.90%  '
$ws.Range('D13').Value = '1.838.73'
$ws.Range('E13').Value = '  +3.02%  '
$ws.Range('D14').Value = '1.608.91'
$ws.Range('E14').Value = '  +3.01%  '
$ws.Range('D15').Value = '29.507.02'
$ws.Range('E15').Value = '  +2.95%  '
$scratch.Value = '0.535'
$scratch.Copy()
$ws.Range('D16').PasteSpecial(-4163)
$ws.Range('E16').Value = '  +3.88%  '
$ws.Range('E17').Value = '  +2.11%  '
$scratch.Value = '63.44'
$scratch.Copy()
$ws.Range('D18').PasteSpecial(-4163)
$ws.Range('E18').Value = '  +3.62%  '
$scratch.Value = '241.04'
$scratch.Copy()
$ws.Range('D19').PasteSpecial(-4163)
$ws.Range('E19').Value = '  +5.67%  '
$ws.Range('E20').Value = '  +4.12%  '
$ws.Range('D21').Value = '0.0₃0689'
$ws.Range('E21').Value = '  +1.79%  '
$ws.Range('E22').Value = '  +0.21%  '
$ws.Range('E23').Value = '  +2.46%  '
$ws.Range('E24').Value = '  +2.17%  '
$ws.Range('E25').Value = '  +0.48%  '
$scratch.Value = '154.49'
$scratch.Copy()
$ws.Range('D26').PasteSpecial(-4163)
$ws.Range('E26').Value = '  +2.21%  '
$ws.Range('E27').Value = '  +4.88%  '
$ws.Range('E28').Value = '  +3.58%  '
$ws.Range('E29').Value = '  +2.63%  '
$ws.Range('E30').Value = '  +0.20%  '
$ws.Range('E31').Value = '  +2.69%  '
$ws.Range('E32').Value = '  +1.20%  '
$ws.Range('E33').Value = '  +1.75%  '
$ws.Range('E34').Value = '  +4.04%  '
$ws.Range('D35').Value = '1.414.69'
$ws.Range('E35').Value = '  +1.67%  '
$ws.Range('E36').Value = '  +0.93%  '
$ws.Range('E37').Value = '  +4.29%  '
$ws.Range('E38').Value = '  +4.78%  '
$ws.Range('E39').Value = '  +0.44%  '
$ws.Range('E40').Value = '  +2.76%  '
$ws.Range('E41').Value = '  +3.63%  '
$ws.Range('E42').Value = '  +1.24%  '
$ws.Range('E43').Value = '  +6.43%  '
$ws.Range('E44').Value = '  +3.54%  '
$ws.Range('E45').Value = '  +0.22%  '
$scratch.Value = '52.46'
$scratch.Copy()
$ws.Range('D46').PasteSpecial(-4163)
$ws.Range('E46').Value = '  +21.28%  '
$scratch.Value = '65.82'
$scratch.Copy()
$ws.Range('D47').PasteSpecial(-4163)
$ws.Range('E47').Value = '  +2.82%  '
$ws.Range('E48').Value = '  +0.96%  '
$ws.Range('D49').Value = '1.748.93'
$ws.Range('E49').Value = '  +3.17%  '
$scratch.Value = '0.857'
$scratch.Copy()
$ws.Range('D50').PasteSpecial(-4163)
$ws.Range('E50').Value = '  -1.37%  '
$scratch.Value = '86.62'
$scratch.Copy()
$ws.Range('D51').PasteSpecial(-4163)
$ws.Range('E51').Value = '  +1.93%  '

# Clean up the scratch cell so it does not leave stray content/formatting behind.
$scratch.Clear()
$excel.CutCopyMode = $false
